$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title paragraph: style Title -> Heading1; text update
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Style = "Heading1"
$p1.Range.Text = "RESIDENTIAL LEASE AGREEMENT"

# ---------------------------------------------------------------------
# 2. Intro sentence
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = 'This Lease Agreement ("Agreement") is entered into on January 1, 2023, by and between:'

# ---------------------------------------------------------------------
# 3. LESSOR line
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = 'LESSOR: UXIN LIMITED ("Landlord")'

# ---------------------------------------------------------------------
# 4. LESSEE line
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Text = 'LESSEE: GLORYFIN INTERNATIONAL GROUP HOLDING COMPANY LIMITED ("Tenant")'

# ---------------------------------------------------------------------
# 5. PROPERTY line -> intro text, then a new "Beijing" paragraph after it
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Text = "PROPERTY: The Landlord hereby leases to the Tenant the residential property located at:"

$p5.Range.InsertParagraphAfter()
$pBeijing = $d.Paragraphs.Item(6)
$pBeijing.Style = "Normal"
$pBeijing.Range.Text = "Beijing"

# ---------------------------------------------------------------------
# 6. "1. TERM OF LEASE" becomes a Heading2
# ---------------------------------------------------------------------
$pTerm = $d.Paragraphs.Item(7)
$pTerm.Style = "Heading2"

# ---------------------------------------------------------------------
# 7. Term-of-lease body text update
# ---------------------------------------------------------------------
$pTermBody = $d.Paragraphs.Item(8)
$pTermBody.Range.Text = "The term of this lease shall commence on January 1, 2023 and shall terminate on December 31, 2023. This Agreement shall be considered a fixed-term lease."

# ---------------------------------------------------------------------
# 8. "2. RENT" becomes a Heading2
# ---------------------------------------------------------------------
$pRent = $d.Paragraphs.Item(9)
$pRent.Style = "Heading2"

# ---------------------------------------------------------------------
# 9. Rent body text update
# ---------------------------------------------------------------------
$pRentBody = $d.Paragraphs.Item(10)
$pRentBody.Range.Text = "The Tenant agrees to pay the Landlord a monthly rent of `$50,000. Rent is due on the 1st day of each month. If rent is not received by the 5th day of the month, a late fee of `$50.00 will be assessed."

# ---------------------------------------------------------------------
# 10. Append new numbered sections 3-16 (Heading2 + Normal body pairs)
# ---------------------------------------------------------------------
$pairs = @(
    @{Heading = '3. SECURITY DEPOSIT'; Body = 'Upon execution of this Agreement, Tenant shall deposit with Landlord the sum of $[SECURITY_DEPOSIT_AMOUNT] as a security deposit. This deposit shall be held by the Landlord as security for the faithful performance by the Tenant of all terms, covenants, and conditions of this Agreement.'},
    @{Heading = '4. USE OF PREMISES'; Body = 'The premises shall be used and occupied by the Tenant exclusively as a private residential dwelling. No part of the premises shall be used for any business or commercial purpose.'},
    @{Heading = '5. UTILITIES'; Body = 'The Tenant shall be responsible for payment of all utilities and services, including but not limited to electricity, gas, water, sewer, trash removal, cable, and internet services.'},
    @{Heading = '6. MAINTENANCE AND REPAIRS'; Body = 'The Tenant shall maintain the premises in good, clean, and sanitary condition. The Tenant shall promptly notify the Landlord of any damage, defects, or dangerous conditions on the premises.'},
    @{Heading = '7. ALTERATIONS'; Body = 'The Tenant shall not make any alterations, additions, or improvements to the premises without the prior written consent of the Landlord.'},
    @{Heading = '8. PETS'; Body = 'No pets shall be allowed on the premises without the prior written consent of the Landlord. If permission is granted, an additional pet deposit may be required.'},
    @{Heading = '9. SUBLETTING'; Body = 'The Tenant shall not sublet the premises or any part thereof, nor assign this lease, without the prior written consent of the Landlord.'},
    @{Heading = '10. RIGHT OF ENTRY'; Body = 'The Landlord reserves the right to enter the premises with reasonable notice (except in cases of emergency) for the purpose of inspection, repairs, or showing the property to prospective tenants or buyers.'},
    @{Heading = '11. TERMINATION'; Body = 'Upon termination of this lease, the Tenant shall surrender the premises in as good condition as when received, reasonable wear and tear excepted. The Tenant shall remove all personal property and clean the premises thoroughly.'},
    @{Heading = '12. DEFAULT'; Body = 'If the Tenant fails to pay rent when due or breaches any other term of this Agreement, the Landlord may terminate this lease with proper notice as required by law.'},
    @{Heading = '13. GOVERNING LAW'; Body = 'This Agreement shall be governed by the laws of the state in which the property is located.'},
    @{Heading = '14. ENTIRE AGREEMENT'; Body = 'This Agreement constitutes the entire agreement between the parties and supersedes all prior negotiations, representations, or agreements, whether written or oral.'},
    @{Heading = '15. SEVERABILITY'; Body = 'If any provision of this Agreement is held to be invalid or unenforceable, the remaining provisions shall continue in full force and effect.'},
    @{Heading = '16. ADDITIONAL TERMS AND CONDITIONS'; Body = '[None specified]'}
)

$last = $d.Paragraphs.Item($d.Paragraphs.Count)
foreach ($pair in $pairs) {
    $last.Range.InsertParagraphAfter()
    $last = $d.Paragraphs.Item($d.Paragraphs.Count)
    $last.Style = "Heading2"
    $last.Range.Text = $pair.Heading

    $last.Range.InsertParagraphAfter()
    $last = $d.Paragraphs.Item($d.Paragraphs.Count)
    $last.Style = "Normal"
    $last.Range.Text = $pair.Body
}

# ---------------------------------------------------------------------
# 11. SIGNATURES heading
# ---------------------------------------------------------------------
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Style = "Heading2"
$last.Range.Text = "SIGNATURES"

# ---------------------------------------------------------------------
# 12. LANDLORD signature block (line breaks, not paragraph breaks)
# ---------------------------------------------------------------------
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Style = "Normal"
$last.Range.Text = "LANDLORD:" + [char]11 + `
    "Signature: _________________________________        Date: _____________" + [char]11 + `
    "Print Name: UXIN LIMITED"

# ---------------------------------------------------------------------
# 13. TENANT signature block
# ---------------------------------------------------------------------
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Style = "Normal"
$last.Range.Text = "TENANT:" + [char]11 + `
    "Signature: _________________________________        Date: _____________" + [char]11 + `
    "Print Name: GLORYFIN INTERNATIONAL GROUP HOLDING COMPANY LIMITED"

# ---------------------------------------------------------------------
# 14. ACKNOWLEDGMENT paragraph
# ---------------------------------------------------------------------
$last.Range.InsertParagraphAfter()
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Style = "Normal"
$last.Range.Text = "ACKNOWLEDGMENT" + [char]11 + `
    "By signing above, both parties acknowledge that they have read, understood, and agree to be bound by all terms and conditions of this Lease Agreement."

Write-Output "edit complete"
